$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: give brand-new cells the correct default column style (A/B/C => s=1/2/3)
# by copying format from a fully-populated template row (row 10: A10=s1, B10=s2, C10=s3)
$ws.Range("A10").Copy($ws.Range("A13"))
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))
$ws.Range("B10").Copy($ws.Range("B23"))
$ws.Range("C10").Copy($ws.Range("C23"))

# --- Step 2: set final cell values
$ws.Cells.Item(10,2).Value = "5840897 - Clodoaldo Saron"
$ws.Cells.Item(10,3).Value = "5840897 - Clodoaldo Saron"
$ws.Cells.Item(13,1).Value = "Programa resumido:"
$ws.Cells.Item(13,2).Value = "Semestral"
$ws.Cells.Item(13,3).Value = "Semestral"
$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Cells.Item(14,2).Value = "Fundamentals of rheology. Processing of polymers: raw materials, machines and molds."
$ws.Cells.Item(14,3).Value = "Fundamentals of rheology. Processing of polymers: raw materials, machines and molds."
$ws.Cells.Item(15,1).Value = "Programa:"
$ws.Cells.Item(15,2).Value = "01/01/2020"
$ws.Cells.Item(15,3).Value = "01/01/2020"
$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Cells.Item(16,2).Value = "1. Introduction to rheology. Types of flow. 2. Hookeian solids and Newtonian fluids. 3. Newtonian and non-Newtonian fluids. 4. Viscoelasticity. Extensional viscosity. Normal stress differences. Variables that affect the viscosity of polymers. 5. Importance of Rheology in the processing of polymers. Flows used to characterize materials: drag flow, flow due to pressure difference and flow in ducts. 6. Extrusion of polymers: equipment, threads, dies and applications. 7. Injection of polymers: equipment, molds, operation control, correction of problems and applications. 8. Other thermoplastic processing techniques: blowing, pressing, thermoforming, calendering, spinning, rotomoulding. 9. Blends and Polymer Composites: obtainable forms, miscibility, compatibility and applications. 10. Processing techniques for thermosetting polymers: manual molding, spray molding, pultrusion, filament winding, pressing, etc."
$ws.Cells.Item(16,3).Value = "1. Introduction to rheology. Types of flow. 2. Hookeian solids and Newtonian fluids. 3. Newtonian and non-Newtonian fluids. 4. Viscoelasticity. Extensional viscosity. Normal stress differences. Variables that affect the viscosity of polymers. 5. Importance of Rheology in the processing of polymers. Flows used to characterize materials: drag flow, flow due to pressure difference and flow in ducts. 6. Extrusion of polymers: equipment, threads, dies and applications. 7. Injection of polymers: equipment, molds, operation control, correction of problems and applications. 8. Other thermoplastic processing techniques: blowing, pressing, thermoforming, calendering, spinning, rotomoulding. 9. Blends and Polymer Composites: obtainable forms, miscibility, compatibility and applications. 10. Processing techniques for thermosetting polymers: manual molding, spray molding, pultrusion, filament winding, pressing, etc."
$ws.Cells.Item(17,1).Value = "Avaliação:"
$ws.Cells.Item(18,1).Value = "Método:"
$ws.Cells.Item(18,2).Value = "5840897 - Clodoaldo Saron"
$ws.Cells.Item(18,3).Value = "5840897 - Clodoaldo Saron"
$ws.Cells.Item(19,1).Value = "Critério:"
$ws.Cells.Item(20,1).Value = "Norma de recuperação:"
$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Cells.Item(22,1).Value = "Requisitos:"
$ws.Cells.Item(23,2).Value = "LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)`n"
$ws.Cells.Item(23,3).Value = "LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)`n"

# --- Step 3: clear cells that should become empty
$ws.Cells.Item(17,2).ClearContents()
$ws.Cells.Item(17,3).ClearContents()
$ws.Cells.Item(22,2).ClearContents()
$ws.Cells.Item(22,3).ClearContents()
$ws.Cells.Item(23,1).ClearContents()

# --- Step 4: row height adjustments
$ws.Rows.Item(13).RowHeight = 60.0
$ws.Rows.Item(15).RowHeight = 120.0
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60.0
$ws.Rows.Item(21).RowHeight = 120.0
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30.0

# --- Step 5: remove row 24 (no longer part of the used range)
$ws.Rows.Item(24).Delete()
